$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 2.9
$ws.Cells.Item(2, 7).Value = 3.45
$ws.Cells.Item(2, 8).Value = 2.28
$ws.Cells.Item(2, 9).Value = 2.58
$ws.Cells.Item(2, 10).Value = 3.35
$ws.Cells.Item(2, 11).Value = 3.9
$ws.Cells.Item(2, 12).Value = 1.37
$ws.Cells.Item(2, 13).Value = 1.05
$ws.Cells.Item(2, 14).Value = 4.5
$ws.Cells.Item(2, 15).Value = 1.24
$ws.Cells.Item(2, 16).Value = 2.16
$ws.Cells.Item(2, 17).Value = 1.81
$ws.Cells.Item(2, 18).Value = 1.46
$ws.Cells.Item(2, 19).Value = 2.9
$ws.Cells.Item(2, 20).Value = 1.56
$ws.Cells.Item(2, 21).Value = 2.38
$ws.Cells.Item(2, 22).Value = 1.64
$ws.Cells.Item(2, 23).Value = 1.41
$ws.Cells.Item(2, 24).Value = 20
$ws.Cells.Item(2, 25).Value = 15
$ws.Cells.Item(2, 26).Value = 18
$ws.Cells.Item(2, 27).Value = 36
$ws.Cells.Item(2, 28).Value = 18.5
$ws.Cells.Item(2, 29).Value = 9.199999999999999
$ws.Cells.Item(2, 30).Value = 12.5
$ws.Cells.Item(2, 31).Value = 26
$ws.Cells.Item(2, 32).Value = 25
$ws.Cells.Item(2, 33).Value = 14.5
$ws.Cells.Item(2, 34).Value = 16.5
$ws.Cells.Item(2, 36).Value = 60
$ws.Cells.Item(2, 37).Value = 36
$ws.Cells.Item(2, 38).Value = 42
$ws.Cells.Item(2, 39).Value = 110
$ws.Cells.Item(2, 40).Value = 25
$ws.Cells.Item(2, 41).Value = 15.5
# Row 3
$ws.Cells.Item(3, 6).Value = 7.2
$ws.Cells.Item(3, 7).Value = 9
$ws.Cells.Item(3, 8).Value = 1.4
$ws.Cells.Item(3, 9).Value = 1.46
$ws.Cells.Item(3, 10).Value = 5.1
$ws.Cells.Item(3, 12).Value = 1.28
$ws.Cells.Item(3, 14).Value = 5.3
$ws.Cells.Item(3, 15).Value = 1.18
$ws.Cells.Item(3, 16).Value = 2.56
$ws.Cells.Item(3, 17).Value = 1.53
$ws.Cells.Item(3, 18).Value = 1.64
$ws.Cells.Item(3, 19).Value = 2.38
$ws.Cells.Item(3, 20).Value = 1.86
$ws.Cells.Item(3, 21).Value = 2
$ws.Cells.Item(3, 22).Value = 3.05
$ws.Cells.Item(3, 23).Value = 1.13
$ws.Cells.Item(3, 24).Value = 44
$ws.Cells.Item(3, 25).Value = 22
$ws.Cells.Item(3, 26).Value = 21
$ws.Cells.Item(3, 27).Value = 26
$ws.Cells.Item(3, 29).Value = 24
$ws.Cells.Item(3, 31).Value = 970
$ws.Cells.Item(3, 33).Value = 80
$ws.Cells.Item(3, 41).Value = 5.6
# Row 4
$ws.Cells.Item(4, 12).Value = 1.4
$ws.Cells.Item(4, 15).Value = 1.32
$ws.Cells.Item(4, 18).Value = 1.34
$ws.Cells.Item(4, 19).Value = 3.05
# Row 5
$ws.Cells.Item(5, 8).Value = 4.6
$ws.Cells.Item(5, 15).Value = 1.4
$ws.Cells.Item(5, 16).Value = 1.74
$ws.Cells.Item(5, 17).Value = 2.1
$ws.Cells.Item(5, 19).Value = 4.1
$ws.Cells.Item(5, 24).Value = 28
$ws.Cells.Item(5, 28).Value = 29
$ws.Cells.Item(5, 32).Value = 40
$ws.Cells.Item(5, 33).Value = 40
$ws.Cells.Item(5, 40).Value = 65
# Row 6
$ws.Cells.Item(6, 11).Value = 4.1
$ws.Cells.Item(6, 12).Value = 1.42
$ws.Cells.Item(6, 15).Value = 1.31
$ws.Cells.Item(6, 17).Value = 1.94
$ws.Cells.Item(6, 20).Value = 1.73
$ws.Cells.Item(6, 21).Value = 2.02
$ws.Cells.Item(6, 22).Value = 1.36
$ws.Cells.Item(6, 24).Value = 14.5
$ws.Cells.Item(6, 25).Value = 14
$ws.Cells.Item(6, 28).Value = 10.5
$ws.Cells.Item(6, 29).Value = 8.6
# Row 7
$ws.Cells.Item(7, 10).Value = 3.9
$ws.Cells.Item(7, 12).Value = 1.39
$ws.Cells.Item(7, 13).Value = 1.07
$ws.Cells.Item(7, 18).Value = 1.34
$ws.Cells.Item(7, 19).Value = 3.15
# Row 8
$ws.Cells.Item(8, 6).Value = 1.84
$ws.Cells.Item(8, 7).Value = 1.94
$ws.Cells.Item(8, 8).Value = 4.9
$ws.Cells.Item(8, 9).Value = 5.4
$ws.Cells.Item(8, 10).Value = 3.45
$ws.Cells.Item(8, 11).Value = 3.9
$ws.Cells.Item(8, 12).Value = 1.5
$ws.Cells.Item(8, 14).Value = 3.05
$ws.Cells.Item(8, 15).Value = 1.42
$ws.Cells.Item(8, 16).Value = 1.68
$ws.Cells.Item(8, 17).Value = 2.28
$ws.Cells.Item(8, 18).Value = 1.25
$ws.Cells.Item(8, 19).Value = 4.3
$ws.Cells.Item(8, 20).Value = 2.02
$ws.Cells.Item(8, 21).Value = 1.79
$ws.Cells.Item(8, 22).Value = 1.23
$ws.Cells.Item(8, 23).Value = 2.06
$ws.Cells.Item(8, 25).Value = 15.5
$ws.Cells.Item(8, 26).Value = 110
$ws.Cells.Item(8, 27).Value = 1000
$ws.Cells.Item(8, 28).Value = 7.4
$ws.Cells.Item(8, 29).Value = 8.4
$ws.Cells.Item(8, 30).Value = 22
$ws.Cells.Item(8, 32).Value = 11
$ws.Cells.Item(8, 34).Value = 25
$ws.Cells.Item(8, 36).Value = 22
$ws.Cells.Item(8, 37).Value = 24
$ws.Cells.Item(8, 40).Value = 19
# Row 9
$ws.Cells.Item(9, 7).Value = 2.54
$ws.Cells.Item(9, 9).Value = 3.25
$ws.Cells.Item(9, 10).Value = 3.4
$ws.Cells.Item(9, 14).Value = 3.7
$ws.Cells.Item(9, 17).Value = 2
$ws.Cells.Item(9, 19).Value = 3.5
$ws.Cells.Item(9, 20).Value = 1.72
$ws.Cells.Item(9, 22).Value = 1.45
$ws.Cells.Item(9, 23).Value = 1.64
$ws.Cells.Item(9, 31).Value = 90
$ws.Cells.Item(9, 32).Value = 27
$ws.Cells.Item(9, 33).Value = 17.5
$ws.Cells.Item(9, 34).Value = 25
$ws.Cells.Item(9, 35).Value = 170
$ws.Cells.Item(9, 36).Value = 140
$ws.Cells.Item(9, 37).Value = 80
$ws.Cells.Item(9, 38).Value = 170
$ws.Cells.Item(9, 39).Value = 330
$ws.Cells.Item(9, 40).Value = 44
$ws.Cells.Item(9, 41).Value = 600
# Row 10
$ws.Cells.Item(10, 6).Value = 1.54
$ws.Cells.Item(10, 7).Value = 2.64
$ws.Cells.Item(10, 8).Value = 2.86
$ws.Cells.Item(10, 11).Value = 4.9
$ws.Cells.Item(10, 16).Value = 1.52
$ws.Cells.Item(10, 17).Value = 1.55
$ws.Cells.Item(10, 18).Value = 1.19
$ws.Cells.Item(10, 19).Value = 1.55
$ws.Cells.Item(10, 22).Value = 1.25
$ws.Cells.Item(10, 23).Value = 1.62
$ws.Cells.Item(10, 24).Value = 30
# Row 11
$ws.Cells.Item(11, 6).Value = 1.61
$ws.Cells.Item(11, 7).Value = 1.63
$ws.Cells.Item(11, 10).Value = 4.7
$ws.Cells.Item(11, 11).Value = 4.9
$ws.Cells.Item(11, 12).Value = 1.29
$ws.Cells.Item(11, 14).Value = 6.4
$ws.Cells.Item(11, 15).Value = 1.17
$ws.Cells.Item(11, 16).Value = 2.78
$ws.Cells.Item(11, 17).Value = 1.52
$ws.Cells.Item(11, 18).Value = 1.71
$ws.Cells.Item(11, 19).Value = 2.3
$ws.Cells.Item(11, 20).Value = 1.63
$ws.Cells.Item(11, 21).Value = 2.48
$ws.Cells.Item(11, 23).Value = 2.58
$ws.Cells.Item(11, 24).Value = 75
$ws.Cells.Item(11, 25).Value = 75
$ws.Cells.Item(11, 26).Value = 370
$ws.Cells.Item(11, 27).Value = 150
$ws.Cells.Item(11, 28).Value = 13.5
$ws.Cells.Item(11, 29).Value = 11
$ws.Cells.Item(11, 31).Value = 65
$ws.Cells.Item(11, 33).Value = 10.5
$ws.Cells.Item(11, 34).Value = 17.5
$ws.Cells.Item(11, 36).Value = 16.5
$ws.Cells.Item(11, 37).Value = 14.5
$ws.Cells.Item(11, 38).Value = 25
$ws.Cells.Item(11, 39).Value = 75
$ws.Cells.Item(11, 40).Value = 6.2
$ws.Cells.Item(11, 41).Value = 120
# Row 12
$ws.Cells.Item(12, 6).Value = 2.98
$ws.Cells.Item(12, 7).Value = 3.1
$ws.Cells.Item(12, 8).Value = 2.48
$ws.Cells.Item(12, 9).Value = 2.54
$ws.Cells.Item(12, 10).Value = 3.6
$ws.Cells.Item(12, 12).Value = 1.38
$ws.Cells.Item(12, 14).Value = 4.5
$ws.Cells.Item(12, 15).Value = 1.26
$ws.Cells.Item(12, 16).Value = 2.18
$ws.Cells.Item(12, 17).Value = 1.8
$ws.Cells.Item(12, 18).Value = 1.48
$ws.Cells.Item(12, 19).Value = 3
$ws.Cells.Item(12, 20).Value = 1.68
$ws.Cells.Item(12, 21).Value = 2.44
$ws.Cells.Item(12, 22).Value = 1.64
$ws.Cells.Item(12, 23).Value = 1.47
$ws.Cells.Item(12, 24).Value = 18
$ws.Cells.Item(12, 25).Value = 13.5
$ws.Cells.Item(12, 28).Value = 15.5
$ws.Cells.Item(12, 29).Value = 8.6
$ws.Cells.Item(12, 30).Value = 12.5
$ws.Cells.Item(12, 32).Value = 23
$ws.Cells.Item(12, 33).Value = 13.5
$ws.Cells.Item(12, 34).Value = 15
$ws.Cells.Item(12, 35).Value = 75
$ws.Cells.Item(12, 36).Value = 55
$ws.Cells.Item(12, 37).Value = 32
$ws.Cells.Item(12, 38).Value = 150
$ws.Cells.Item(12, 39).Value = 200
$ws.Cells.Item(12, 40).Value = 25
$ws.Cells.Item(12, 41).Value = 19
# Row 13
$ws.Cells.Item(13, 6).Value = 3.55
$ws.Cells.Item(13, 7).Value = 3.7
$ws.Cells.Item(13, 8).Value = 2.04
$ws.Cells.Item(13, 9).Value = 2.06
$ws.Cells.Item(13, 12).Value = 1.31
$ws.Cells.Item(13, 14).Value = 5.5
$ws.Cells.Item(13, 16).Value = 2.54
$ws.Cells.Item(13, 17).Value = 1.61
$ws.Cells.Item(13, 18).Value = 1.61
$ws.Cells.Item(13, 19).Value = 2.52
$ws.Cells.Item(13, 20).Value = 1.58
$ws.Cells.Item(13, 21).Value = 2.54
$ws.Cells.Item(13, 22).Value = 1.94
$ws.Cells.Item(13, 23).Value = 1.37
$ws.Cells.Item(13, 24).Value = 32
$ws.Cells.Item(13, 25).Value = 13.5
$ws.Cells.Item(13, 26).Value = 15.5
$ws.Cells.Item(13, 27).Value = 48
$ws.Cells.Item(13, 28).Value = 21
$ws.Cells.Item(13, 29).Value = 10
$ws.Cells.Item(13, 30).Value = 10.5
$ws.Cells.Item(13, 31).Value = 18.5
$ws.Cells.Item(13, 33).Value = 16
$ws.Cells.Item(13, 36).Value = 70
$ws.Cells.Item(13, 37).Value = 36
$ws.Cells.Item(13, 38).Value = 40
$ws.Cells.Item(13, 40).Value = 26
$ws.Cells.Item(13, 41).Value = 10.5
# Row 14
$ws.Cells.Item(14, 8).Value = 1.56
$ws.Cells.Item(14, 9).Value = 1.58
$ws.Cells.Item(14, 10).Value = 4.5
$ws.Cells.Item(14, 11).Value = 4.7
$ws.Cells.Item(14, 14).Value = 4.8
$ws.Cells.Item(14, 15).Value = 1.23
$ws.Cells.Item(14, 16).Value = 2.26
$ws.Cells.Item(14, 17).Value = 1.75
$ws.Cells.Item(14, 18).Value = 1.49
$ws.Cells.Item(14, 19).Value = 2.84
$ws.Cells.Item(14, 20).Value = 1.85
$ws.Cells.Item(14, 21).Value = 2.12
$ws.Cells.Item(14, 22).Value = 2.72
$ws.Cells.Item(14, 24).Value = 19.5
$ws.Cells.Item(14, 28).Value = 26
$ws.Cells.Item(14, 29).Value = 10
$ws.Cells.Item(14, 30).Value = 9.6
$ws.Cells.Item(14, 31).Value = 15
$ws.Cells.Item(14, 32).Value = 60
$ws.Cells.Item(14, 33).Value = 25
$ws.Cells.Item(14, 34).Value = 20
$ws.Cells.Item(14, 36).Value = 210
$ws.Cells.Item(14, 38).Value = 85
$ws.Cells.Item(14, 40).Value = 1000
